# Update cryptos list values (price and 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.845.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.581.51"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.478"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.04"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.802.21"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.576.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.81%  "

$ws.Range("E14").Value = "  -2.88%  "

$ws.Range("E15").Value = "  -2.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.817.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("E18").Value = "  -3.30%  "

$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.33"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.13%  "

$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("E29").Value = "  -3.27%  "

$ws.Range("E30").Value = "  -5.69%  "

$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.12"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("E33").Value = "  -2.62%  "

$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.097.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.83%  "

$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.38%  "

$ws.Range("E39").Value = "  -2.03%  "

$ws.Range("E40").Value = "  -3.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.777"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -8.06%  "

$ws.Range("E42").Value = "  +7.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.31%  "

$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.716.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("E51").Value = "  -0.27%  "
